# Auto-applies updated crypto price/volume values to sheet1 (D/E columns, rows 2-51).
# Values are written as literal text (matching the source inlineStr cells), so for
# any Price value that Excel would otherwise auto-parse as a number we temporarily
# force the cell to Text format, then restore the default ("Normal") style so no
# stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2: '35.307.68' -> '35.346.32'; E2: '  +0.50%  ' -> '  +0.60%  '
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.346.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.60%  "

# Row 3: D3: '1.881.51' -> '1.881.67'; E3: '  -1.06%  ' -> '  -1.04%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.881.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.04%  "

# Row 4: E4: '  -0.65%  ' -> '  -0.68%  '
$ws.Range("E4").Value = "  -0.68%  "

# Row 5: D5: '245.13' -> '245.03'; E5: '  -3.22%  ' -> '  -3.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.27%  "

# Row 6: D6: '0.689' -> '0.687'; E6: '  -0.71%  ' -> '  -0.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.687"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.80%  "

# Row 7: E7: '  -0.73%  ' -> '  -0.77%  '
$ws.Range("E7").Value = "  -0.77%  "

# Row 8: D8: '43.49' -> '43.58'; E8: '  +5.47%  ' -> '  +5.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.58"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.18%  "

# Row 9: E9: '  -0.97%  ' -> '  -0.86%  '
$ws.Range("E9").Value = "  -0.86%  "

# Row 10: D10: '53.33' -> '53.39'; E10: '  +0.99%  ' -> '  +1.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.19%  "

# Row 11: E11: '  -1.39%  ' -> '  -1.35%  '
$ws.Range("E11").Value = "  -1.35%  "

# Row 12: E12: '  -1.00%  ' -> '  -1.12%  '
$ws.Range("E12").Value = "  -1.12%  "

# Row 13: D13: '13.33' -> '13.35'; E13: '  +2.17%  ' -> '  +2.67%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.35"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.67%  "

# Row 14: D14: '2.155.30' -> '2.155.10'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.155.10"
$ws.Range("D14").Style = "Normal"

# Row 15: D15: '0.756' -> '0.757'; E15: '  +2.98%  ' -> '  +3.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.757"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.20%  "

# Row 16: D16: '4.89' -> '4.90'; E16: '  -1.68%  ' -> '  -1.67%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.67%  "

# Row 17: D17: '1.893.44' -> '1.896.88'; E17: '  -0.40%  ' -> '  -0.24%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.896.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.24%  "

# Row 18: D18: '35.411.67' -> '35.413.14'; E18: '  +0.80%  ' -> '  +0.79%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "35.413.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.79%  "

# Row 19: D19: '72.89' -> '72.92'; E19: '  -0.85%  ' -> '  -0.87%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.87%  "

# Row 20: D20: '0.0₃0821' -> '0.0₃0820'; E20: '  -1.41%  ' -> '  -1.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0820"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.36%  "

# Row 21: D21: '244.24' -> '244.20'; E21: '  +0.57%  ' -> '  +0.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "244.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.35%  "

# Row 22: E22: '  -1.39%  ' -> '  -1.47%  '
$ws.Range("E22").Value = "  -1.47%  "

# Row 23: E23: '  -2.29%  ' -> '  -2.15%  '
$ws.Range("E23").Value = "  -2.15%  "

# Row 24: D24: '2.65' -> '2.66'; E24: '  +8.79%  ' -> '  +9.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +9.31%  "

# Row 25: E25: '  -0.67%  ' -> '  -0.68%  '
$ws.Range("E25").Value = "  -0.68%  "

# Row 26: D26: '2.13' -> '2.14'; E26: '  -6.66%  ' -> '  -6.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.57%  "

# Row 27: D27: '165.36' -> '165.31'; E27: '  -0.98%  ' -> '  -0.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.97%  "

# Row 28: D28: '8.50' -> '8.51'; E28: '  -0.55%  ' -> '  -0.60%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.60%  "

# Row 29: D29: '18.28' -> '18.29'; E29: '  -1.20%  ' -> '  -1.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.24%  "

# Row 30: E30: '  -1.97%  ' -> '  -2.02%  '
$ws.Range("E30").Value = "  -2.02%  "

# Row 31: D31: '4.128.46' -> '4.128.45'; E31: '  +0.01%  ' -> '  +0.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.128.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.00%  "

# Row 32: D32: '1.69' -> '1.70'; E32: '  +7.58%  ' -> '  +7.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.33%  "

# Row 33: E33: '  -1.30%  ' -> '  -1.32%  '
$ws.Range("E33").Value = "  -1.32%  "

# Row 34: D34: '0.0585' -> '0.0587'; E34: '  -4.07%  ' -> '  -3.57%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0587"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.57%  "

# Row 35: D35: '1.90' -> '1.89'; E35: '  -7.08%  ' -> '  -7.76%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.89"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.76%  "

# Row 36: E36: '  -1.65%  ' -> '  -1.63%  '
$ws.Range("E36").Value = "  -1.63%  "

# Row 37: E37: '  -0.72%  ' -> '  -0.71%  '
$ws.Range("E37").Value = "  -0.71%  "

# Row 38: D38: '0.842' -> '0.843'; E38: '  -1.33%  ' -> '  -1.59%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.843"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.59%  "

# Row 39: E39: '  -2.50%  ' -> '  -2.60%  '
$ws.Range("E39").Value = "  -2.60%  "

# Row 40: D40: '0.0696' -> '0.0697'; E40: '  +7.08%  ' -> '  +7.18%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0697"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.18%  "

# Row 41: D41: '17.18' -> '17.28'; E41: '  -1.00%  ' -> '  -0.74%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.74%  "

# Row 42: D42: '0.0218' -> '0.0217'; E42: '  +1.19%  ' -> '  +0.97%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0217"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.97%  "

# Row 43: D43: '96.07' -> '96.01'; E43: '  -6.88%  ' -> '  -6.71%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "96.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.71%  "

# Row 44: E44: '  -2.48%  ' -> '  -2.31%  '
$ws.Range("E44").Value = "  -2.31%  "

# Row 45: D45: '1.303.18' -> '1.304.66'; E45: '  -1.24%  ' -> '  -1.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.304.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.18%  "

# Row 46: E46: '  -4.54%  ' -> '  -4.26%  '
$ws.Range("E46").Value = "  -4.26%  "

# Row 47: D47: '0.0795' -> '0.0798'; E47: '  +6.63%  ' -> '  +6.96%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0798"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.96%  "

# Row 48: D48: '2.38' -> '2.37'; E48: '  -2.06%  ' -> '  -2.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.23%  "

# Row 49: D49: '2.73' -> '2.72'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.72"
$ws.Range("D49").Style = "Normal"

# Row 50: D50: '12.14' -> '12.15'; E50: '  +1.59%  ' -> '  +0.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "12.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.33%  "

# Row 51: E51: '  -5.48%  ' -> '  -5.49%  '
$ws.Range("E51").Value = "  -5.49%  "

